$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Insert a new row at position 27 (pushes old rows 27..69 down to 28..70)
# ---------------------------------------------------------------------------
$ws.Rows.Item(27).Insert()

# ---------------------------------------------------------------------------
# 2. Fix up existing cell content that changed (rows above / at the insert)
# ---------------------------------------------------------------------------

# Row 12 : add a new bug report + state "Corrigé"
$ws.Range("C12").Value2 = "incohérence : le formulaire est nommé ""absence rapide"", et lorsqu'une absence est ajoutée, ""Retard ajouté"" est affiché. (+ faute de frappe)"
$ws.Range("D12").Value2 = "Corrigé"

# Row 16 : reworded bug description
$ws.Range("C16").Value2 = "bug détecté : le contact et le stagiaire sont liés bizarrement, ça ne tient pas trop la route. La correction apporte une gestion plus souple."

# Row 24 : remove the old "Abandonné" state value
$ws.Range("D24").ClearContents()

# Row 25 : shrink the row height, replace the long bug text with the shorter one,
# and change the state from "En cours" to "Corrigé"
$ws.Rows.Item(25).RowHeight = 30
$ws.Range("C25").Value2 = "bug détecté : on ne peut pas ajouter ni modifier de contact - idem pour les entreprises… `nbug détecté : Même quand la requête ne réussit pas, la popup affiche ""suppression effectuée"". "
$ws.Range("D25").Value2 = "Corrigé"

# ---------------------------------------------------------------------------
# 3. Populate the newly inserted row 27
# ---------------------------------------------------------------------------
$ws.Range("B27").Value2 = "A compléter avec les opérations sur les contacts, les entreprises….."
$ws.Range("C27").Value2 = "bug détecté : Lors de la modification d'un contact, on peut lui enlever ses 2 numéros de téléphone"
$ws.Range("D27").Value2 = "Corrigé"

# ---------------------------------------------------------------------------
# 4. Sheet view: scroll position / selected cell
# ---------------------------------------------------------------------------
$ws.Application.ActiveWindow.ScrollRow = 34
$ws.Range("C53").Select()

Write-Host "Done"
